$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Insert a new blank row at 402, pushing current 402..446 down to 403..447.
$ws.Rows.Item(402).Insert()

# 2. Fix up formatting of the brand-new row 402: copy the (already-correct,
#    shifted-down) formats from row 403 across every table column. This
#    reuses the existing style indices instead of minting near-duplicate
#    ones, matching what Excel itself does when expanding a table row.
$ws.Range("A403:K403").Copy()
$ws.Range("A402:K402").PasteSpecial(-4122)

# 3. The REMARKS/date column (K) uses a distinct date-formatted style (the
#    same one already used a couple of rows up, e.g. K401) for the two new
#    leave-card rows; copy that style onto K402 and K403.
$ws.Range("K401").Copy()
$ws.Range("K402").PasteSpecial(-4122)
$ws.Range("K403").PasteSpecial(-4122)

# 4. Restore the calculated-column formula in the new row 402's EARNED(2)
#    cell (row insert leaves it blank).
$ws.Range("G402").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),`"`",Table1[[#This Row],[EARNED]])"

# 5. Populate the new leave entry (row 402) - SL(1-0-0), 1 day, dated 9/29/2023.
$ws.Range("B402").Value = "SL(1-0-0)"
$ws.Range("H402").Value = 1
$ws.Range("K402").Value = 45198

# 6. Row 403 (formerly row 402) becomes the second SL(1-0-0) entry: 1.25
#    days earned, 1 day charged, dated 10/31/2023.
$ws.Range("B403").Value = "SL(1-0-0)"
$ws.Range("C403").Value = 1.25
$ws.Range("H403").Value = 1
$ws.Range("K403").Value = 45230

# 7. Row 404 (formerly row 403) becomes an SL(2-0-0) entry, 2 days charged,
#    remarks holding the free-text date range instead of a single date.
$ws.Range("B404").Value = "SL(2-0-0)"
$ws.Range("H404").Value = 2
$ws.Range("K404").Value = "11/13,14/2023"

# 8. Grow the table definition to include the newly appended last row.
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A8:K447"))

# 9. Make sure the calculated-column formula text on the (shifted) final
#    table row stays in the canonical structured-reference form.
$ws.Range("G447").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),`"`",Table1[[#This Row],[EARNED]])"

# 10. Recalculate so every dependent cell (BALANCE columns, CONVERTION!A7,
#     etc.) carries a fresh cached value.
$excel.CalculateFullRebuild()

# 11. Leave the on-screen selection near the rows that were just edited,
#     mirroring where a person would naturally be after typing this in.
$ws.Application.Goto($ws.Range("K404"), $true)
$ws.Range("K404").Select()
